$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-14"

# Update the header label for the current year column (shared string)
$ws.Range("I1").Value = "2022 (through 08-14)"

# Update the August row value (row 9) for column I
$ws.Range("I9").Value = 80

# Update the Total row value (row 14) for column I
$ws.Range("I14").Value = 1050
